$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "matrices" scores (column F) - small re-simulation deltas
$ws.Range("F2").Value = 14.36820170540361
$ws.Range("F3").Value = 13.10445131131576
$ws.Range("F4").Value = 8.165085299637123
$ws.Range("F5").Value = 7.471687554102642
$ws.Range("F6").Value = 6.394477316427598
$ws.Range("F7").Value = 6.210208915507454
$ws.Range("F8").Value = 5.464785281559637

# Rows 9 and 10 swap ranking order: the two individuals (prolificid/name/race)
# trade places because their recomputed scores cross over.
$ws.Range("B9").Value = 33
$ws.Range("C9").Value = "60b322994d0b901954690036"
$ws.Range("D9").Value = "Brennan"
$ws.Range("F9").Value = 5.311456126118004
$ws.Range("G9").Value = "White"

$ws.Range("B10").Value = 32
$ws.Range("C10").Value = "60bf9943e4e04642d4634ecc"
$ws.Range("D10").Value = "Jamarii"
$ws.Range("F10").Value = 5.011805588421218
$ws.Range("G10").Value = "Black or African American"

$ws.Range("F11").Value = 3.497412773125043
$ws.Range("F12").Value = 1.210064642988239
$ws.Range("F13").Value = 0.4333488792121737
